$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8; this shifts old rows 8..115 down to 9..116
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new data entry.
# Most fields mirror the row that was previously in row 8 (now at row 9);
# only Fecha, Volumen, Precio minimo/maximo/promedio and Precio $/Kg differ.
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C8").Value = "Los Lagos"
$ws.Range("D8").Value = 44699
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 27000
$ws.Range("L8").Value = 27000
$ws.Range("M8").Value = 27000
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región de La Araucanía"
$ws.Range("P8").Value = 1080
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
